$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# The underlying report lists pharmacy items sorted alphabetically.
# Four new items were added to the report (HEBTA C.M.D TAB, SORAL 30 MG
# 30CAPS, VONACIDAN 20 MG 20 F.C.TABS., and "معجون اسنان فلورو بالكولا"),
# which pushes the existing rows (and the totals / footer rows below
# them) down by four rows, and the grand total is recalculated.
# -----------------------------------------------------------------------

# Step 1: insert 4 new blank rows right before the totals row (row 19)
$ws.Range("A19:A22").EntireRow.Insert()

# Step 2: clone the formatting (styles) of an existing item row (row 18)
# into the newly inserted rows
$ws.Range("A18:N18").Copy()
$ws.Range("A19:N22").PasteSpecial(-4122)

# Step 3: re-create the merged cell layout for each new row
$ws.Range("B19:G19").Merge()
$ws.Range("H19:K19").Merge()
$ws.Range("L19:M19").Merge()
$ws.Range("B20:G20").Merge()
$ws.Range("H20:K20").Merge()
$ws.Range("L20:M20").Merge()
$ws.Range("B21:G21").Merge()
$ws.Range("H21:K21").Merge()
$ws.Range("L21:M21").Merge()
$ws.Range("B22:G22").Merge()
$ws.Range("H22:K22").Merge()
$ws.Range("L22:M22").Merge()

# Step 4: row heights for the new rows
$ws.Rows.Item(19).RowHeight = 24.75
$ws.Rows.Item(20).RowHeight = 25.5
$ws.Rows.Item(21).RowHeight = 24.75
$ws.Rows.Item(22).RowHeight = 25.5

# Step 5: footer row height changes slightly now that it moved to row 24
$ws.Rows.Item(24).RowHeight = 17.25

# Step 6: sequence numbers (column A) for the 4 new rows
$ws.Range("A19").Value = 16
$ws.Range("A20").Value = 17
$ws.Range("A21").Value = 18
$ws.Range("A22").Value = 19

# Step 7: (re)write the full, re-sorted item list into rows 10-22.
# columns: B = item name, H = current balance ratio, L = sale price,
# N = number of transactions ratio
$ws.Range("B10").Value = "HEBTA C.M.D TAB"
$ws.Range("H10").Value = "0:0"
$ws.Range("L10").Value = 144
$ws.Range("N10").Value = "1:0"

$ws.Range("B11").Value = "IVYROSPAN SYRUP 100 ML"
$ws.Range("H11").Value = "0:0"
$ws.Range("L11").Value = 55
$ws.Range("N11").Value = "1:0"

$ws.Range("B12").Value = "KETOLAC 30MG/2ML 5 AMP. FOR I.M./I.V. INF."
$ws.Range("H12").Value = "3:2"
$ws.Range("L12").Value = 12
$ws.Range("N12").Value = "0:0"

$ws.Range("B13").Value = "MAXILASE 200 CEIP UNIT/ML SYRUP 100ML"
$ws.Range("H13").Value = "4:0"
$ws.Range("L13").Value = 57
$ws.Range("N13").Value = "1:0"

$ws.Range("B14").Value = "METAPSIN 10 F.C. TABS."
$ws.Range("H14").Value = "1:0"
$ws.Range("L14").Value = 110
$ws.Range("N14").Value = "1:0"

$ws.Range("B15").Value = "PICOLAX 0.75% ORAL DROPS 15 ML"
$ws.Range("H15").Value = "3:0"
$ws.Range("L15").Value = 23.04
$ws.Range("N15").Value = "1:0"

$ws.Range("B16").Value = "SORAL 30 MG 30CAPS"
$ws.Range("H16").Value = "0:0"
$ws.Range("L16").Value = 86
$ws.Range("N16").Value = "0:3"

$ws.Range("B17").Value = "VOLTAREN 75MG/3ML 3 AMP."
$ws.Range("H17").Value = "3:2"
$ws.Range("L17").Value = 17
$ws.Range("N17").Value = "0:0"

$ws.Range("B18").Value = "VONACIDAN 20 MG 20 F.C.TABS."
$ws.Range("H18").Value = "0:0"
$ws.Range("L18").Value = 192
$ws.Range("N18").Value = "1:0"

$ws.Range("B19").Value = "WATER FOR INJECTION AMP. 5 ML"
$ws.Range("H19").Value = "7789:0"
$ws.Range("L19").Value = 5
$ws.Range("N19").Value = "2:0"

$ws.Range("B20").Value = "بلاستر مترسيلك 2.5 سم"
$ws.Range("H20").Value = "36:0"
$ws.Range("L20").Value = 25
$ws.Range("N20").Value = "1:0"

$ws.Range("B21").Value = "سرنجات 3 سم"
$ws.Range("H21").Value = "-1:0"
$ws.Range("L21").Value = 16
$ws.Range("N21").Value = "8:0"

$ws.Range("B22").Value = "معجون اسنان فلورو بالكولا"
$ws.Range("H22").Value = "3:0"
$ws.Range("L22").Value = 30
$ws.Range("N22").Value = "1:0"

# Step 8: update the grand total (now on row 23 after the insert)
$ws.Range("K23").Value = 1218.04
